$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Parameter labels in column H ---
# "% of body mass" -> "% of body mass (daily ration)"
$ws.Range("H26").Value = "% of body mass (daily ration)"
$ws.Range("H30").Value = "% of body mass (daily ration)"
$ws.Range("H34").Value = "% of body mass (daily ration)"
$ws.Range("H38").Value = "% of body mass (daily ration)"

# "ADMR" -> "Assimilation rate"
$ws.Range("H27").Value = "Assimilation rate"
$ws.Range("H31").Value = "Assimilation rate"
$ws.Range("H35").Value = "Assimilation rate"
$ws.Range("H39").Value = "Assimilation rate"

# "Assimilation rate" -> "Average Daily Metabolic Rate (kJ)"
$ws.Range("H28").Value = "Average Daily Metabolic Rate (kJ)"
$ws.Range("H32").Value = "Average Daily Metabolic Rate (kJ)"
$ws.Range("H36").Value = "Average Daily Metabolic Rate (kJ)"
$ws.Range("H40").Value = "Average Daily Metabolic Rate (kJ)"

# --- Swap the B:G numeric blocks between the ADMR row and the Assimilation-rate row ---
# (the row labels stayed put, but which metric each row actually holds changed)
foreach ($col in 2..7) {
    $tmp = $ws.Cells.Item(27, $col).Value()
    $ws.Cells.Item(27, $col).Value = $ws.Cells.Item(28, $col).Value()
    $ws.Cells.Item(28, $col).Value = $tmp
}
foreach ($col in 2..7) {
    $tmp = $ws.Cells.Item(31, $col).Value()
    $ws.Cells.Item(31, $col).Value = $ws.Cells.Item(32, $col).Value()
    $ws.Cells.Item(32, $col).Value = $tmp
}
foreach ($col in 2..7) {
    $tmp = $ws.Cells.Item(35, $col).Value()
    $ws.Cells.Item(35, $col).Value = $ws.Cells.Item(36, $col).Value()
    $ws.Cells.Item(36, $col).Value = $tmp
}
foreach ($col in 2..7) {
    $tmp = $ws.Cells.Item(39, $col).Value()
    $ws.Cells.Item(39, $col).Value = $ws.Cells.Item(40, $col).Value()
    $ws.Cells.Item(40, $col).Value = $tmp
}

# --- Updated numeric results (re-run of the underlying model) ---
$ws.Cells.Item(22, 2).Value = 14.0361928983377
$ws.Cells.Item(22, 3).Value = 41.1870222447896
$ws.Cells.Item(22, 4).Value = 67.8596324642761
$ws.Cells.Item(22, 5).Value = 67.8614401171181
$ws.Cells.Item(22, 6).Value = 94.275170507753
$ws.Cells.Item(22, 7).Value = 125.936954519803

$ws.Cells.Item(23, 2).Value = 4.86455732322552
$ws.Cells.Item(23, 3).Value = 14.1276755508576
$ws.Cells.Item(23, 4).Value = 23.2219415296271
$ws.Cells.Item(23, 5).Value = 23.21854697563
$ws.Cells.Item(23, 6).Value = 32.263068737367
$ws.Cells.Item(23, 7).Value = 45.5170354499567

$ws.Cells.Item(24, 2).Value = 2.65575157778109
$ws.Cells.Item(24, 3).Value = 6.90817659179424
$ws.Cells.Item(24, 4).Value = 11.4150196133521
$ws.Cells.Item(24, 5).Value = 11.4137593981445
$ws.Cells.Item(24, 6).Value = 15.8984011910272
$ws.Cells.Item(24, 7).Value = 20.6635091441203

$ws.Cells.Item(25, 2).Value = 3.48686200365584
$ws.Cells.Item(25, 3).Value = 9.86833419581917
$ws.Cells.Item(25, 4).Value = 16.2392247715859
$ws.Cells.Item(25, 5).Value = 16.2513431554639
$ws.Cells.Item(25, 6).Value = 22.5886556034896
$ws.Cells.Item(25, 7).Value = 30.0102835019049

$ws.Cells.Item(42, 2).Value = 137.555780271871
$ws.Cells.Item(42, 3).Value = 372.382082685988
$ws.Cells.Item(42, 4).Value = 838.036399849542
$ws.Cells.Item(42, 5).Value = 785.90375444921
$ws.Cells.Item(42, 6).Value = 1606.6733330904
$ws.Cells.Item(42, 7).Value = 4780.24966983913

$ws.Cells.Item(43, 2).Value = 27.205246559231
$ws.Cells.Item(43, 3).Value = 98.3825781560992
$ws.Cells.Item(43, 4).Value = 226.022479175145
$ws.Cells.Item(43, 5).Value = 210.142845883117
$ws.Cells.Item(43, 6).Value = 446.085147127326
$ws.Cells.Item(43, 7).Value = 1426.71880776252

$ws.Cells.Item(44, 2).Value = 8.49111604926644
$ws.Cells.Item(44, 3).Value = 38.682812294605
$ws.Cells.Item(44, 4).Value = 87.3069414379729
$ws.Cells.Item(44, 5).Value = 81.2109972800295
$ws.Cells.Item(44, 6).Value = 170.980370623951
$ws.Cells.Item(44, 7).Value = 568.861810121852

$ws.Cells.Item(45, 2).Value = 15.0272345185261
$ws.Cells.Item(45, 3).Value = 49.8254742761078
$ws.Cells.Item(45, 4).Value = 111.749974905978
$ws.Cells.Item(45, 5).Value = 104.03866354813
$ws.Cells.Item(45, 6).Value = 218.805544920763
$ws.Cells.Item(45, 7).Value = 572.921456279388

$ws.Cells.Item(46, 2).Value = 111.827476605863
$ws.Cells.Item(46, 3).Value = 292.828769887502
$ws.Cells.Item(46, 4).Value = 670.491869650502
$ws.Cells.Item(46, 5).Value = 626.494098574606
$ws.Cells.Item(46, 6).Value = 1300.53328516584
$ws.Cells.Item(46, 7).Value = 4130.66591475121

$ws.Cells.Item(47, 2).Value = 20.4905754748067
$ws.Cells.Item(47, 3).Value = 77.5807158517612
$ws.Cells.Item(47, 4).Value = 180.722562045428
$ws.Cells.Item(47, 5).Value = 167.66623106851
$ws.Cells.Item(47, 6).Value = 360.395961262783
$ws.Cells.Item(47, 7).Value = 1272.27303818768

$ws.Cells.Item(48, 2).Value = 6.4061200810256
$ws.Cells.Item(48, 3).Value = 30.4547822719123
$ws.Cells.Item(48, 4).Value = 69.8499235605624
$ws.Cells.Item(48, 5).Value = 64.7783597128404
$ws.Cells.Item(48, 6).Value = 138.638060820435
$ws.Cells.Item(48, 7).Value = 484.100055227918

$ws.Cells.Item(49, 2).Value = 11.7075493787266
$ws.Cells.Item(49, 3).Value = 39.2040944835421
$ws.Cells.Item(49, 4).Value = 89.4361575926233
$ws.Cells.Item(49, 5).Value = 83.0298349242873
$ws.Cells.Item(49, 6).Value = 176.987297479381
$ws.Cells.Item(49, 7).Value = 487.869453168769

Write-Output "edit.ps1 applied"
